$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Paragraph 1: blank spacer paragraph appended right after the final
# existing paragraph of the document (inherits the plain, non-bold,
# left-aligned Times New Roman / 12pt formatting already in effect at
# the end of the story).
# ------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.Text = "`r"

# ------------------------------------------------------------------
# Paragraph 2: centred, bold heading "Ένατη συνάντηση. (2/5/2021)"
# ------------------------------------------------------------------
$headRange = $d.Content
$headRange.Collapse(0)
$headRange.Text = "Ένατη συνάντηση. (2/5/2021)`r"

$count = $d.Paragraphs.Count
$headPara = $d.Paragraphs.Item($count - 1)
$headParaRange = $headPara.Range
$headParaRange.Font.Name = "Times New Roman"
$headParaRange.Font.NameBi = "Times New Roman"
$headParaRange.Font.Size = 12
$headParaRange.Font.SizeBi = 12
$headParaRange.Font.Bold = 1
$headParaRange.ParagraphFormat.Alignment = 1

# ------------------------------------------------------------------
# Paragraph 3: blank paragraph, still centred + bold (matches the
# paragraph mark formatting under the heading). Typed as a temporary
# placeholder character so the bold/centre formatting actually takes,
# then the placeholder character is removed again, leaving a clean
# empty paragraph whose mark keeps the centred/bold formatting.
# ------------------------------------------------------------------
$blankRange = $d.Content
$blankRange.Collapse(0)
$blankRange.Text = "X`r"

$count = $d.Paragraphs.Count
$blankPara = $d.Paragraphs.Item($count - 1)
$blankParaRange = $blankPara.Range
$blankParaRange.Font.Name = "Times New Roman"
$blankParaRange.Font.NameBi = "Times New Roman"
$blankParaRange.Font.Size = 12
$blankParaRange.Font.SizeBi = 12
$blankParaRange.Font.Bold = 1
$blankParaRange.ParagraphFormat.Alignment = 1

$placeholder = $d.Range($blankPara.Range.Start, $blankPara.Range.Start + 1)
$placeholder.Delete()

# ------------------------------------------------------------------
# Paragraph 4: normal (left aligned, non-bold) body paragraph with the
# write-up of the ninth meeting.
# ------------------------------------------------------------------
$bodyRange = $d.Content
$bodyRange.Collapse(0)
$bodyRange.Text = "Με την ένατη συνάντηση και την εφαρμογή να έχει φτάσει σχεδόν στο τέλος της, συζητήθηκαν απαιτήσεις και πως θα γίνει η υλοποίηση αυτών από τους προγραμματιστές για την καλύτερη δυνατή λειτουργικότητα από τον μελλοντικό χρήστη. Έπειτα της συνεννόησης των προγραμματιστών και πελατών συζητήθηκε η γενική πρόοδος των μελών καθώς και τι αλλαγές θα πρέπει να γίνουν έτσι ώστε η εργασία να έχει το καλύτερο δυνατόν παρουσιαστικό σύμφωνα με το έγγραφο που μας παρείχε ο καθηγητής."

$count = $d.Paragraphs.Count
$bodyPara = $d.Paragraphs.Item($count)
$bodyParaRange = $bodyPara.Range
$bodyParaRange.Font.Name = "Times New Roman"
$bodyParaRange.Font.NameBi = "Times New Roman"
$bodyParaRange.Font.Size = 12
$bodyParaRange.Font.SizeBi = 12
$bodyParaRange.Font.Bold = 0
$bodyParaRange.ParagraphFormat.Alignment = 0
